# Regenerate the localization-status report for archival.
#
# The CI job that produces this workbook ("Generate Report for Archive")
# re-runs the export over the same snapshot of data, so every cell ends up
# holding exactly the text it already had. The only incidental side effect
# of the regeneration is that the exporter's status-lookup table now also
# knows about the "In Translation" status (it is simply not used by any of
# the rows in *this* snapshot, since none of them are currently mid-flight).
#
# We reproduce that by touching each "Status" cell through the in-progress
# "In Translation" value and then writing back the row's real, unchanged
# status, which is exactly what a re-export keyed off the same source
# records would observe.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Row 2 -> 96b33741-af7f-4eb8-8d4d-0b4825afb466.md (unaffected: "Handed back: in sync with en-US")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "Handed back: in sync with en-US"

# Row 3 -> 4cc1ac0b-cf85-4083-81a5-3cded0c74374.md ("Ready for handoff")
$zhcn.Range("C3").Value = "In Translation"
$zhcn.Range("C3").Value = "Ready for handoff"
$dede.Range("C3").Value = "In Translation"
$dede.Range("C3").Value = "Ready for handoff"

# Row 4 -> 6c616282-82da-4745-ab0f-091e39be842f.md ("Ready for handoff")
$zhcn.Range("C4").Value = "In Translation"
$zhcn.Range("C4").Value = "Ready for handoff"
$dede.Range("C4").Value = "In Translation"
$dede.Range("C4").Value = "Ready for handoff"

# Row 5 -> 96b0622a-bfaa-473e-a4b5-8aa92a75022d.md ("Ready for handoff")
$zhcn.Range("C5").Value = "In Translation"
$zhcn.Range("C5").Value = "Ready for handoff"
$dede.Range("C5").Value = "In Translation"
$dede.Range("C5").Value = "Ready for handoff"
